$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16 (pushes existing rows 16-23 down to 17-24)
$ws.Rows.Item(16).Insert()

# Populate the new row 16 with the "Steering Rack protection" part data
$ws.Cells.Item(16, 3).Value = "Steering Rack protection"
$ws.Cells.Item(16, 4).Value = "m"
$ws.Cells.Item(16, 5).Value = "To protect the steering rack."
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = "ST_03003"

# Copy style of the adjacent data row (row 15) to the new row for consistent formatting
$ws.Range("A15:G15").Copy()
$ws.Range("A16:G16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Match row height of a standard (non customHeight) data row
$ws.Rows.Item(16).RowHeight = $ws.Rows.Item(14).RowHeight

# Update selection to match final state
$ws.Range("J38").Select()
